# Word COM-interop script: append the "SAMPLE TEST PLAN (on canvas)" content
# (plus the three follow-on paragraphs describing the sample test plan outline,
# comments and CS3300 template) after the existing body text, replacing the
# single trailing empty paragraph.
#
# NOTE: Range.InsertXML() on this host silently clobbers the *entire* document
# body instead of inserting at the target range when asked to insert 3+
# paragraphs in one call, so the new content is inserted in two smaller
# batches (<=2 paragraphs each) to stay well clear of that limit.

$d = $word.ActiveDocument

# --- Batch 1: the bold "SAMPLE TEST PLAN (on canvas):" heading paragraph and
#     the big "You will have multiple test plans..." outline paragraph.
$target = $d.Paragraphs.Last.Range
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>SAMPLE TEST PLAN (on canvas):</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>You will have multiple test plans. Test plans are normally numbered.</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>A test plan should reflect ‘HOW’ you plan to test, along with ‘WHAT’ you plan to test and</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>“WHEN” you plan to execute the test.</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>Sample Test Plan Outline</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>1. Project Title: ABC</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>2. Test Plan ###</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>3. Test Objectives: Unit, Integration, or System test</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>4. Test Approach:</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>a. First test with correct input data. (ensures execution)</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>b. Second test with out of boundary input date. (defect examination)</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>c. Third test functionality. (ensure requirements are successfully met)</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>d. Fourth test performance, stress, dependability (</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>5. Manual or Automated Test:</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>6. Test tools: NONE</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>7. Test environment: Joe’s laptop</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>8. Test criteria: Requirement ID numbers from the requirements spreadsheet/ requirements</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>traceability matrix. Functional or Non-Functional Requirements. Which requirements</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>are you testing?</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>9. Test schedule: When do you plan to test this. Sunday, March 9, 2025</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">10. Test team (list by name): Note: the unit developer </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>can not</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> be the unit tester.</w:t>
      </w:r>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$target.InsertXML($xml1)

# --- Batch 2: the "COMMENTS:" paragraph and the "CS3300 PROJECT NAME..." paragraph.
$target = $d.Paragraphs.Last.Range
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:r>
        <w:t>COMMENTS:</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">1, Functional testing objective – Ensure software works as it should. Validate user </w:t>
      </w:r>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>workflows,</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">data processing, and verify input/output parameters. </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Check for out of</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> boundary alerts.</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>2. Performance testing objectives – ensure software is efficient and can handle various loads</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>(throughput and scalability)</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>3. Security testing objectives – Uncover security flaws</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>4. Usability testing – Verify user workflow</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>CS3300 PROJECT NAME TEST PLAN ###</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>DATE</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>A Test Strategy or a Detailed Test Plan has EACH input step identified along with the expected</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>output. No button is pushed arbitrarily. These are very large documents. One such plan</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>normally takes 2 to 3 days to produce, but for large complex systems it may take 7 to 10 days.</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/>
        </w:rPr>
        <w:t>➢</w:t>
      </w:r>
    </w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$target.InsertXML($xml2)

# The four new paragraphs above were inserted just *before* the document's
# original trailing empty paragraph, which is now redundant (it used to be
# the only paragraph after the "...response" text). Remove it so the new
# content flows directly into the section properties, matching the target.
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$secondLastPara = $d.Paragraphs.Item($paraCount - 1)
$trailingMark = $d.Range($secondLastPara.Range.End - 1, $lastPara.Range.End)
$trailingMark.Delete()
